# Applies the 2023-12-07 cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # These D/E columns are stored as text (inline strings) in the workbook,
    # even though many values look numeric (e.g. "1.00", "233.52", "2.698.44").
    # Prefixing with a leading apostrophe forces Excel to keep them as text
    # instead of silently coercing them into real numbers.
    $ws.Range($cellRef).Value = "`'" + $text
}

$ws.Range("D2").Value = "43.731.94"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "2.349.60"
$ws.Range("E3").Value = "  +3.52%  "
Set-TextCell "D4" "1.00"
$ws.Range("E4").Value = "  +0.03%  "
Set-TextCell "D5" "233.52"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  +2.11%  "
Set-TextCell "D7" "65.79"
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextCell "D9" "0.454"
$ws.Range("E9").Value = "  +1.30%  "
Set-TextCell "D10" "0.0972"
$ws.Range("E10").Value = "  -4.47%  "
Set-TextCell "D11" "56.75"
$ws.Range("E11").Value = "  -0.47%  "
Set-TextCell "D12" "26.88"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "2.698.44"
$ws.Range("E13").Value = "  +3.51%  "
$ws.Range("E14").Value = "  -0.96%  "
Set-TextCell "D15" "15.44"
$ws.Range("E15").Value = "  -1.50%  "
Set-TextCell "D16" "6.18"
$ws.Range("E16").Value = "  -0.72%  "
Set-TextCell "D17" "0.854"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "2.347.57"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").Value = "43.675.81"
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").Value = "  -2.82%  "
Set-TextCell "D21" "74.03"
$ws.Range("E21").Value = "  +0.49%  "
Set-TextCell "D22" "6.26"
$ws.Range("E22").Value = "  +2.94%  "
Set-TextCell "D23" "249.85"
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D24" "1.00"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("B25").Value = "WEMIXToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D25" "3.80"
$ws.Range("E25").Value = "  +14.60%  "
Set-TextCell "D26" "2.43"
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("E27").Value = "  +0.16%  "
Set-TextCell "D28" "9.95"
$ws.Range("E28").Value = "  -0.89%  "
Set-TextCell "D29" "22.37"
$ws.Range("E29").Value = "  +7.33%  "
Set-TextCell "D30" "175.03"
$ws.Range("E30").Value = "  +1.85%  "
Set-TextCell "D31" "1.44"
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("E32").Value = "  -5.43%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("E34").Value = "  +4.37%  "
$ws.Range("E35").Value = "  -1.79%  "
Set-TextCell "D36" "4.99"
$ws.Range("E36").Value = "  +1.65%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D37" "3.75"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D38" "2.45"
$ws.Range("E38").Value = "  +5.77%  "
Set-TextCell "D39" "6.59"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  +11.63%  "
Set-TextCell "D42" "1.00"
$ws.Range("E42").Value = "  +0.04%  "
Set-TextCell "D43" "17.97"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("E44").Value = "  +10.57%  "
Set-TextCell "D45" "99.57"
$ws.Range("E45").Value = "  +1.43%  "
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("E47").Value = "  -0.29%  "
Set-TextCell "D48" "4.37"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("B49").Value = "Celestia"
$ws.Range("C49").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell "D49" "10.09"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "1.448.78"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("E51").Value = "  +0.10%  "
